# Daily attendance processing - 2025-10-09 22:18:37
# Swap the order of the first two "Recorded By" entries (column G) for the
# session rows that were re-recorded / re-synced, so that the most recent
# recorder is listed first (e.g. "System, X" -> "X, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,6,11,12,13,14,15,29,30,33,38,39,40,41,42,56,57,58,60,65,66,67,68,69,86,89,90,93,95,112,115,116,119,121,138,141,142,145,147)

foreach ($row in $rows) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $current = [string]$cell.Value2
    $parts = $current -split ', '
    if ($parts.Length -ge 2) {
        $first = $parts[0]
        $second = $parts[1]
        $parts[0] = $second
        $parts[1] = $first
        $cell.Value = [string]::Join(', ', $parts)
    }
}
